# "Fruta / hortaliza, semanal" - weekly update: a new price observation
# for Haba (Vega Modelo de Temuco) is inserted as a new row 42, pushing
# the existing historical rows (old 42-55) down by one (new rows 43-56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 42 (shifts rows 42..55 down to 43..56,
# and extends the used range / dimension to row 56 automatically).
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with this week's data.
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44523
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 100112026
$ws.Range("G42").Value = "Haba"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 20
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 8000
$ws.Range("M42").Value = 8000
$ws.Range("N42").Value = '$/saco 25 kilos'
$ws.Range("O42").Value = "Región de La Araucanía"
$ws.Range("P42").Value = 320
$ws.Range("Q42").Value = 25
$ws.Range("R42").Value = "Hortaliza"
